# Rename the "_old"/"_new" column-header suffixes to the respective input
# file's format version ("_FV2210" / "_FV2304"), add a real Excel Table
# (ListObject) over the whole used range with an AutoFilter, and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Header renames -----------------------------------------------------
# Columns A:J were suffixed "_old" -> now "_FV2210"
$headersFV2210 = @(
    "Segmentname_FV2210",
    "Segmentgruppe_FV2210",
    "Segment_FV2210",
    "Datenelement_FV2210",
    "Segment ID_FV2210",
    "Code_FV2210",
    "Qualifier_FV2210",
    "Beschreibung_FV2210",
    "Bedingungsausdruck_FV2210",
    "Bedingung_FV2210"
)
for ($i = 0; $i -lt $headersFV2210.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2210[$i]
}

# K1 ("diff") stays untouched.

# Columns L:U were suffixed "_new" -> now "_FV2304"
$headersFV2304 = @(
    "Segmentname_FV2304",
    "Segmentgruppe_FV2304",
    "Segment_FV2304",
    "Datenelement_FV2304",
    "Segment ID_FV2304",
    "Code_FV2304",
    "Qualifier_FV2304",
    "Beschreibung_FV2304",
    "Bedingungsausdruck_FV2304",
    "Bedingung_FV2304"
)
for ($i = 0; $i -lt $headersFV2304.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $headersFV2304[$i]
}

# --- 2. Turn the used range into a real Table (ListObject) -----------------
$usedRange = $ws.UsedRange
$tbl = $ws.ListObjects.Add(1, $usedRange, 0, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row -----------------------------------------------
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
